$wb = $excel.ActiveWorkbook

# Both "展览" and "全部类型" sheets contain the same F2/F3/F5 "想去人数" values
# that need to be updated to the latest scraped counts.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 770
    $ws.Range("F3").Value = 4177
    $ws.Range("F5").Value = 759
}
